$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @([double]"0.9999969141016266", [double]"0.9990169585194482", [double]"0.9999999996986246", [double]"0.9999997859112689", [double]"0.9999999294435367", [double]"2.88055067624043e-06", [double]"0.0009176260715613002", [double]"5.047781821777987e-10", [double]"1.759675952515807e-07", [double]"8.823618671687925e-08", [double]"9.999963621106398e-05", [double]"0.00169721851163615", [double]"0.9999753128130129", [double]"0.001769472617403162", [double]"67.51505815013246", [double]"93.11145047236465")

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
